$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: forecast week starting dates (serial date numbers), shifted forward one week
$row1 = @(45613, 45620, 45627, 45634, 45641, 45648, 45655, 45662, 45669, 45676, 45683)
for ($i = 0; $i -lt $row1.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $row1[$i]
}

# Rows 2-12: updated VCI3M forecast values per location (columns B..L)
$data = @{
    2  = @(71, 67.2, 63.6, 60.3, 57.6, 55.7, 54.7, 54.5, 55.1, 56.3, 57.9)
    3  = @(72.3, 67.5, 62.8, 58.6, 55.1, 52.4, 50.5, 49.5, 49.3, 49.7, 50.4)
    4  = @(66.59999999999999, 65.7, 64.90000000000001, 64.3, 64, 64.2, 64.7, 65.59999999999999, 66.8, 68.09999999999999, 69.40000000000001)
    5  = @(71.3, 68.09999999999999, 65.09999999999999, 62.3, 59.9, 57.9, 56.3, 54.9, 53.6, 52.3, 51)
    6  = @(82.3, 75.7, 68.8, 62.1, 55.9, 50.6, 46.2, 43, 41, 40.2, 40.3)
    7  = @(68.7, 68.2, 68.2, 68.8, 70.2, 72.3, 75.3, 79, 83.09999999999999, 87.40000000000001, 91.59999999999999)
    8  = @(67.7, 66.7, 66.2, 66.5, 67.5, 69.5, 72.5, 76.3, 80.7, 85.40000000000001, 90.09999999999999)
    9  = @(60.1, 55.3, 50.9, 47, 44.2, 42.7, 42.7, 44.1, 47, 51, 55.7)
    10 = @(71.59999999999999, 66.3, 60.9, 55.6, 50.8, 46.5, 43, 40.4, 38.5, 37.5, 37.3)
    11 = @(86.90000000000001, 78, 68.8, 59.9, 52, 45.4, 40.8, 38.1, 37.6, 39.1, 42.2)
    12 = @(62.3, 60.8, 59.3, 58.1, 56.9, 55.9, 55, 54.1, 53.2, 52.3, 51.3)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
